# Appends the new benchmark-run rows (31-38) produced by the latest
# hyperparameter sweep to the results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = 20
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 0.003
$ws.Cells.Item(31, 5).Value = "FGSM"
$ws.Cells.Item(31, 9).Value = "<function relu at 0x128d789d8>"
$ws.Cells.Item(31, 10).Value = 0.9329000115394592
$ws.Cells.Item(31, 11).Value = 0.7416999936103821
$ws.Cells.Item(31, 12).Value = 0.5971999764442444
$ws.Cells.Item(31, 13).Value = 0.2206476628780365
$ws.Cells.Item(31, 14).Value = 0.74456787109375
$ws.Cells.Item(31, 15).Value = 0.7416999936103821
$ws.Cells.Item(31, 16).Value = "logs/results_328.log"
$ws.Cells.Item(31, 17).Value = "weights/model_328.ckpt"
$ws.Cells.Item(31, 18).Value = "tb/328"
$ws.Cells.Item(31, 19).Value = "(12.799518, 16.56418, 25.464622, 25.766693, 14.9766655, 6.904501, 3.7851212)"
$ws.Cells.Item(31, 20).Value = "(309.03058, 13.829748, 16.421766, 12.924315, 13.326067, 11.006118, 10.691159, 13.276878)"

# Row 32
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).Value = 30
$ws.Cells.Item(32, 3).Value = 0.0005
$ws.Cells.Item(32, 4).Value = 0.003
$ws.Cells.Item(32, 5).Value = "Regular"
$ws.Cells.Item(32, 9).Value = "<function relu at 0x113fe39d8>"
$ws.Cells.Item(32, 10).Value = 0.9472000002861023
$ws.Cells.Item(32, 11).Value = 0.3325000107288361
$ws.Cells.Item(32, 12).Value = 0.1923000067472458
$ws.Cells.Item(32, 13).Value = 0.1799724549055099
$ws.Cells.Item(32, 14).Value = 3.20237398147583
$ws.Cells.Item(32, 15).Value = 0.3325000107288361
$ws.Cells.Item(32, 16).Value = "logs/results_331.log"
$ws.Cells.Item(32, 17).Value = "weights/model_331.ckpt"
$ws.Cells.Item(32, 18).Value = "tb/331"
$ws.Cells.Item(32, 19).Value = "(1.8040013, 1.8094578, 1.9554862, 2.5528448, 3.399231, 4.635884, 6.39143)"
$ws.Cells.Item(32, 20).Value = "(28.089184, 4.674623, 3.6114936, 3.2815678, 3.1776936, 2.116904, 2.2179012, 2.6026235)"

# Row 33
$ws.Cells.Item(33, 1).Value = 31
$ws.Cells.Item(33, 2).Value = 30
$ws.Cells.Item(33, 3).Value = 0.0009
$ws.Cells.Item(33, 4).Value = 0.003
$ws.Cells.Item(33, 5).Value = "Regular"
$ws.Cells.Item(33, 9).Value = "<function relu at 0x118f1b9d8>"
$ws.Cells.Item(33, 10).Value = 0.8956000208854675
$ws.Cells.Item(33, 11).Value = 0.3217999935150146
$ws.Cells.Item(33, 12).Value = 0.2044000029563904
$ws.Cells.Item(33, 13).Value = 0.3738516569137573
$ws.Cells.Item(33, 14).Value = 3.476548433303833
$ws.Cells.Item(33, 15).Value = 0.3217999935150146
$ws.Cells.Item(33, 16).Value = "logs/results_332.log"
$ws.Cells.Item(33, 17).Value = "weights/model_332.ckpt"
$ws.Cells.Item(33, 18).Value = "tb/332"
$ws.Cells.Item(33, 19).Value = "(1.4544966, 1.1530415, 1.2943891, 2.1282125, 4.040144, 7.1142826, 6.1251535)"
$ws.Cells.Item(33, 20).Value = "(21.872171, 4.500764, 3.899859, 2.0486226, 2.4281611, 2.2534535, 1.7587844, 2.4474022)"

# Row 34
$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(34, 2).Value = 60
$ws.Cells.Item(34, 3).Value = 0.0015
$ws.Cells.Item(34, 4).Value = 0.0003
$ws.Cells.Item(34, 5).Value = "Regular"
$ws.Cells.Item(34, 9).Value = "<function relu at 0x10b0f09d8>"
$ws.Cells.Item(34, 10).Value = 0.9347000122070312
$ws.Cells.Item(34, 11).Value = 0.2572999894618988
$ws.Cells.Item(34, 12).Value = 0.1509999930858612
$ws.Cells.Item(34, 13).Value = 0.2364677786827087
$ws.Cells.Item(34, 14).Value = 4.247204780578613
$ws.Cells.Item(34, 15).Value = 0.2572999894618988
$ws.Cells.Item(34, 16).Value = "logs/results_333.log"
$ws.Cells.Item(34, 17).Value = "weights/model_333.ckpt"
$ws.Cells.Item(34, 18).Value = "tb/333"
$ws.Cells.Item(34, 19).Value = "(0.4250503, 0.4371529, 0.74407357, 1.2906153, 2.5302472, 4.9194384, 6.6656017)"
$ws.Cells.Item(34, 20).Value = "(6.92537, 5.7589846, 3.67073, 3.039194, 3.1920185, 5.1396704, 3.4247825, 4.5205474)"

# Row 35
$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(35, 2).Value = 60
$ws.Cells.Item(35, 3).Value = 0.000015
$ws.Cells.Item(35, 4).Value = 0.0003
$ws.Cells.Item(35, 5).Value = "Regular"
$ws.Cells.Item(35, 9).Value = "<function relu at 0x1113b09d8>"
$ws.Cells.Item(35, 10).Value = 0.9394999742507935
$ws.Cells.Item(35, 11).Value = 0.09109999984502792
$ws.Cells.Item(35, 12).Value = 0.03319999948143959
$ws.Cells.Item(35, 13).Value = 0.221238300204277
$ws.Cells.Item(35, 14).Value = 6.377280712127686
$ws.Cells.Item(35, 15).Value = 0.09109999984502792
$ws.Cells.Item(35, 16).Value = "logs/results_336.log"
$ws.Cells.Item(35, 17).Value = "weights/model_336.ckpt"
$ws.Cells.Item(35, 18).Value = "tb/336"
$ws.Cells.Item(35, 19).Value = "(2.0724185, 1.3095005, 1.6185311, 1.3509899, 1.5398465, 1.9591715, 4.275175)"
$ws.Cells.Item(35, 20).Value = "(39.100674, 3.319005, 3.0519383, 1.6382135, 2.70516, 2.7627785, 5.8854613, 6.491546)"

# Row 36
$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(36, 2).Value = 30
$ws.Cells.Item(36, 3).Value = 0.000023
$ws.Cells.Item(36, 4).Value = 0.0003
$ws.Cells.Item(36, 5).Value = "Regular"
$ws.Cells.Item(36, 9).Value = "<function relu at 0x1132949d8>"
$ws.Cells.Item(36, 10).Value = 0.9189000129699707
$ws.Cells.Item(36, 11).Value = 0.05869999900460243
$ws.Cells.Item(36, 12).Value = 0.02209999971091747
$ws.Cells.Item(36, 13).Value = 0.3333184719085693
$ws.Cells.Item(36, 14).Value = 6.367059707641602
$ws.Cells.Item(36, 15).Value = 0.05869999900460243
$ws.Cells.Item(36, 16).Value = "logs/results_346.log"
$ws.Cells.Item(36, 17).Value = "weights/model_346.ckpt"
$ws.Cells.Item(36, 18).Value = "tb/346"
$ws.Cells.Item(36, 19).Value = "(1.5257524, 0.7150309, 0.56886697, 0.49205807, 0.57274866, 1.5198021, 4.5895967)"
$ws.Cells.Item(36, 20).Value = "(33.718807, 2.9477801, 2.0959005, 1.7632871, 2.2024982, 4.748903, 8.546004, 6.910356)"

# Row 37
$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(37, 2).Value = 100
$ws.Cells.Item(37, 3).Value = 0.007
$ws.Cells.Item(37, 4).Value = 0.0003
$ws.Cells.Item(37, 5).Value = "Regular"
$ws.Cells.Item(37, 9).Value = "<function relu at 0x1226c79d8>"
$ws.Cells.Item(37, 10).Value = 0.9646999835968018
$ws.Cells.Item(37, 11).Value = 0.1290999948978424
$ws.Cells.Item(37, 12).Value = 0.01889999955892563
$ws.Cells.Item(37, 13).Value = 0.1261469274759293
$ws.Cells.Item(37, 14).Value = 6.858481407165527
$ws.Cells.Item(37, 15).Value = 0.1290999948978424
$ws.Cells.Item(37, 16).Value = "logs/results_353.log"
$ws.Cells.Item(37, 17).Value = "weights/model_353.ckpt"
$ws.Cells.Item(37, 18).Value = "tb/353"
$ws.Cells.Item(37, 19).Value = "(0.6469947, 0.73037505, 0.9719021, 1.3849623, 2.158197, 3.6967816, 6.6338224)"
$ws.Cells.Item(37, 20).Value = "(9.72191, 3.5845523, 3.115291, 2.8329852, 2.7452235, 2.6821659, 2.723758, 2.6171021)"

# Row 38
$ws.Cells.Item(38, 1).Value = 36
$ws.Cells.Item(38, 2).Value = 100
$ws.Cells.Item(38, 3).Value = 0.07
$ws.Cells.Item(38, 4).Value = 0.0003
$ws.Cells.Item(38, 5).Value = "Regular"
$ws.Cells.Item(38, 9).Value = "<function relu at 0x1226499d8>"
$ws.Cells.Item(38, 10).Value = 0.9384999871253967
$ws.Cells.Item(38, 11).Value = 0.3431999981403351
$ws.Cells.Item(38, 12).Value = 0.217399999499321
$ws.Cells.Item(38, 13).Value = 0.2118040025234222
$ws.Cells.Item(38, 14).Value = 2.253601551055908
$ws.Cells.Item(38, 15).Value = 0.3431999981403351
$ws.Cells.Item(38, 16).Value = "logs/results_354.log"
$ws.Cells.Item(38, 17).Value = "weights/model_354.ckpt"
$ws.Cells.Item(38, 18).Value = "tb/354"
$ws.Cells.Item(38, 19).Value = "(0.2376155, 0.2859654, 0.43169716, 0.67097735, 1.0554945, 1.7692584, 2.8368058)"
$ws.Cells.Item(38, 20).Value = "(3.0683854, 1.7650508, 1.7416625, 1.7273158, 1.7367424, 1.746138, 1.7387878, 1.7745363)"

